$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Repayment schedule" sheet: add a new column O (values of 0) mirroring
#    the existing column N for rows 2-14 (row 2/3 stay blank, rows 4-14 get
#    an explicit 0), keeping the same cell style as column N.
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

$wsSchedule.Range("N2:N14").Copy() | Out-Null
$wsSchedule.Range("O2:O14").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

for ($r = 4; $r -le 14; $r++) {
    $wsSchedule.Cells.Item($r, 15).Value = 0
}

# ---------------------------------------------------------------------------
# 2. "Summary" sheet: move the stored selection to D3.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate() | Out-Null
$wsSummary.Range("D3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. "Transactions" sheet: update the loan/transaction ids and move the
#    stored selection to D3, restoring it as the active tab.
# ---------------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Activate() | Out-Null

$wsTransactions.Range("A2").Value = 7078
$wsTransactions.Range("A3").Value = 7077
$wsTransactions.Range("A4").Value = 7076

$wsTransactions.Range("D3").Select() | Out-Null
